$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Refresh the "Fixed" date placeholder text (1/16/2026 -> 1/30/2026)
#    across the slide master and every slide layout (mirrors using the
#    Header & Footer dialog with "Apply to All").
# ------------------------------------------------------------------
$m = $p.SlideMaster

for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $sh = $m.Shapes.Item($j)
    if ($sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = "1/30/2026"
    }
}

for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
    $cl = $m.CustomLayouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $sh = $cl.Shapes.Item($j)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = "1/30/2026"
        }
    }
}

$p.NotesMaster.HeadersFooters.DateAndTime.Text = "1/30/2026"

# ------------------------------------------------------------------
# 2) Update the title on the first slide.
# ------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(4).TextFrame.TextRange.Text = "Logistic Regression"

# ------------------------------------------------------------------
# 3) Nudge "Picture 7" on slide 6 up slightly.
# ------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(5).Top = 273.78914

Write-Output "edits applied"
